$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '25.953.50'
$ws.Range("E2").Value = '  +0.39%  '
$ws.Range("D3").Value = '1.586.78'
$ws.Range("E3").Value = '  +0.05%  '
$ws.Range("D5").Value = "'210.15"
$ws.Range("E5").Value = '  +0.13%  '
$ws.Range("E6").Value = '  -0.17%  '
$ws.Range("D7").Value = "'0.479"
$ws.Range("E7").Value = '  -0.14%  '
$ws.Range("D8").Value = "'0.245"
$ws.Range("E9").Value = '  -1.11%  '
$ws.Range("D10").Value = "'17.84"
$ws.Range("E10").Value = '  -1.16%  '
$ws.Range("D11").Value = "'0.0806"
$ws.Range("E11").Value = '  +2.02%  '
$ws.Range("D12").Value = '1.808.87'
$ws.Range("E12").Value = '  +0.12%  '
$ws.Range("D13").Value = '1.588.91'
$ws.Range("E13").Value = '  +0.19%  '
$ws.Range("E14").Value = '  -1.51%  '
$ws.Range("D15").Value = "'0.508"
$ws.Range("E15").Value = '  -0.14%  '
$ws.Range("D16").Value = '25.938.53'
$ws.Range("E16").Value = '  +0.30%  '
$ws.Range("D17").Value = "'59.90"
$ws.Range("E17").Value = '  +0.09%  '
$ws.Range("D18").Value = '0.0₃0718'
$ws.Range("E18").Value = '  -0.63%  '
$ws.Range("E19").Value = '  -0.13%  '
$ws.Range("D20").Value = "'198.83"
$ws.Range("E20").Value = '  +3.81%  '
$ws.Range("E21").Value = '  +0.54%  '
$ws.Range("D22").Value = "'9.14"
$ws.Range("E22").Value = '  -2.34%  '
$ws.Range("E23").Value = '  +0.46%  '
$ws.Range("E24").Value = '  +8.75%  '
$ws.Range("D25").Value = "'142.33"
$ws.Range("E25").Value = '  +0.23%  '
$ws.Range("E26").Value = '  -0.06%  '
$ws.Range("E27").Value = '  -8.82%  '
$ws.Range("D28").Value = "'14.99"
$ws.Range("E28").Value = '  -0.81%  '
$ws.Range("D29").Value = "'6.43"
$ws.Range("E29").Value = '  -0.37%  '
$ws.Range("E30").Value = '  -0.22%  '
$ws.Range("D31").Value = "'0.0472"
$ws.Range("E31").Value = '  +0.53%  '
$ws.Range("D32").Value = "'3.10"
$ws.Range("E32").Value = '  -0.23%  '
$ws.Range("D33").Value = "'2.92"
$ws.Range("E33").Value = '  -3.60%  '
$ws.Range("B34").Value = 'LidoDAOToken'
$ws.Range("C34").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D34").Value = "'1.46"
$ws.Range("E34").Value = '  -2.41%  '
$ws.Range("B35").Value = 'HuobiToken'
$ws.Range("C35").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D35").Value = "'2.36"
$ws.Range("E35").Value = '  +0.26%  '
$ws.Range("D36").Value = '1.120.43'
$ws.Range("E36").Value = '  +1.72%  '
$ws.Range("E37").Value = '  +7.17%  '
$ws.Range("E38").Value = '  -0.09%  '
$ws.Range("E39").Value = '  -1.54%  '
$ws.Range("D40").Value = "'0.778"
$ws.Range("E40").Value = '  +0.29%  '
$ws.Range("D41").Value = "'0.486"
$ws.Range("E41").Value = '  -3.66%  '
$ws.Range("D42").Value = "'0.777"
$ws.Range("E42").Value = '  -5.30%  '
$ws.Range("D43").Value = '1.719.90'
$ws.Range("E43").Value = '  +0.00%  '
$ws.Range("E44").Value = '  -2.32%  '
$ws.Range("D45").Value = "'91.74"
$ws.Range("E45").Value = '  -2.21%  '
$ws.Range("B46").Value = 'Aave'
$ws.Range("C46").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D46").Value = "'53.06"
$ws.Range("E46").Value = '  -0.30%  '
$ws.Range("B47").Value = 'RenderToken'
$ws.Range("C47").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D47").Value = "'1.47"
$ws.Range("E47").Value = '  -2.75%  '
$ws.Range("E48").Value = '  -1.11%  '
$ws.Range("D49").Value = "'0.407"
$ws.Range("E49").Value = '  +0.08%  '
$ws.Range("E50").Value = '  +0.18%  '
